$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("principles of special and temporal locality", $true, $false, $false, $false, $false, $true, 1, $false, "principles of spatial and temporal locality", 2)
